$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (G14) onto the
# new header cell (H14) before writing its value, so the new column inherits
# the same font/border style used by the rest of the header row.
$ws.Range("G14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H14").Value = "Case Type"

# Nudge the width of the newly introduced columns to match the authored layout.
$ws.Columns.Item(7).ColumnWidth = 12.33
$ws.Columns.Item(8).ColumnWidth = 11.33

# Match the cursor position left after adding the column.
$ws.Range("E14").Select() | Out-Null
